$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 158 (weekly price record), pushing the
# existing rows 158..166 down to 159..167.
$ws.Rows(158).Insert()

# Populate the new row with the same market/product attributes as its
# neighbours, but for the next week's date (2022-01-24 = serial 44585)
# and the same volume/price figures as the (now shifted) row 159.
$ws.Cells.Item(158, 1).Value = 4
$ws.Cells.Item(158, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(158, 3).Value = "Los Lagos"
$ws.Cells.Item(158, 4).Value = 44585
$ws.Cells.Item(158, 5).Value = 10
$ws.Cells.Item(158, 6).Value = 100112039
$ws.Cells.Item(158, 7).Value = "Ciboulette"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 80
$ws.Cells.Item(158, 11).Value = 3000
$ws.Cells.Item(158, 12).Value = 3000
$ws.Cells.Item(158, 13).Value = 3000
$ws.Cells.Item(158, 14).Value = "`$/docena de atados"
$ws.Cells.Item(158, 15).Value = "Región Metropolitana"
$ws.Cells.Item(158, 16).Value = 1000
$ws.Cells.Item(158, 17).Value = 3
$ws.Cells.Item(158, 18).Value = "Hortaliza"
